$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestInputData")

# Previously: A1=Browser, B1=Product ; A2=Chrome, B2=Samsung Galaxy S24 Ultra
# New layout adds a dedicated column per browser (Chrome, Edge) with Y/N support flags,
# keeping the Product column (shifted right) last.
$ws.Range("B:B").Insert() | Out-Null

$ws.Range("A1").Value = "Chrome"
$ws.Range("B1").Value = "Edge"
$ws.Range("C1").Value = "Product"

$ws.Range("A2").Value = "Y"
$ws.Range("B2").Value = "Y"
$ws.Range("C2").Value = "Samsung Galaxy S24 Ultra"

$ws.Range("A2").Select()
